$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11
$ws.Cells.Item(11, 7).Value = 2.02
$ws.Cells.Item(11, 17).Value = 2.32

# Row 20
$ws.Cells.Item(20, 14).Value = 1.75
$ws.Cells.Item(20, 15).Value = 2.05

# Row 25
$ws.Cells.Item(25, 7).Value = 3.95
$ws.Cells.Item(25, 8).Value = 3.25
$ws.Cells.Item(25, 9).Value = 1.88
$ws.Cells.Item(25, 12).Value = 1.38
$ws.Cells.Item(25, 13).Value = 2.6
$ws.Cells.Item(25, 14).Value = 2.1
$ws.Cells.Item(25, 15).Value = 1.57
$ws.Cells.Item(25, 18).Value = 1.93
$ws.Cells.Item(25, 19).Value = 1.7
$ws.Cells.Item(25, 21).Value = 20
$ws.Cells.Item(25, 22).Value = 13.5
$ws.Cells.Item(25, 23).Value = 60
$ws.Cells.Item(25, 24).Value = 40
$ws.Cells.Item(25, 25).Value = 55
$ws.Cells.Item(25, 26).Value = 7.8
$ws.Cells.Item(25, 27).Value = 6.4
$ws.Cells.Item(25, 28).Value = 17.5
$ws.Cells.Item(25, 29).Value = 100
$ws.Cells.Item(25, 30).Value = 1000
$ws.Cells.Item(25, 31).Value = 5.9
$ws.Cells.Item(25, 32).Value = 8
$ws.Cells.Item(25, 34).Value = 15.5
$ws.Cells.Item(25, 35).Value = 17
$ws.Cells.Item(25, 36).Value = 35

# Row 26
$ws.Cells.Item(26, 7).Value = 2.55
$ws.Cells.Item(26, 8).Value = 2.95
$ws.Cells.Item(26, 9).Value = 2.8
$ws.Cells.Item(26, 12).Value = 1.44
$ws.Cells.Item(26, 13).Value = 2.4
$ws.Cells.Item(26, 14).Value = 2.27
$ws.Cells.Item(26, 15).Value = 1.5
$ws.Cells.Item(26, 18).Value = 1.93
$ws.Cells.Item(26, 19).Value = 1.7
$ws.Cells.Item(26, 20).Value = 6.8
$ws.Cells.Item(26, 21).Value = 11.75
$ws.Cells.Item(26, 22).Value = 10
$ws.Cells.Item(26, 23).Value = 28
$ws.Cells.Item(26, 24).Value = 24
$ws.Cells.Item(26, 27).Value = 5.8
$ws.Cells.Item(26, 28).Value = 16.5
$ws.Cells.Item(26, 29).Value = 100
$ws.Cells.Item(26, 31).Value = 6.9
$ws.Cells.Item(26, 32).Value = 12.5
$ws.Cells.Item(26, 33).Value = 10.75
$ws.Cells.Item(26, 34).Value = 35
$ws.Cells.Item(26, 35).Value = 29
$ws.Cells.Item(26, 36).Value = 45

# Row 28
$ws.Cells.Item(28, 7).Value = 1.75
$ws.Cells.Item(28, 9).Value = 4.7
$ws.Cells.Item(28, 20).Value = 7.4
$ws.Cells.Item(28, 21).Value = 8.75
$ws.Cells.Item(28, 24).Value = 13
$ws.Cells.Item(28, 27).Value = 6.6
$ws.Cells.Item(28, 28).Value = 14
$ws.Cells.Item(28, 31).Value = 12.5
$ws.Cells.Item(28, 32).Value = 28
$ws.Cells.Item(28, 35).Value = 50

# Row 37
$ws.Cells.Item(37, 7).Value = 2.72
$ws.Cells.Item(37, 8).Value = 3.55
$ws.Cells.Item(37, 9).Value = 2.32
$ws.Cells.Item(37, 10).Value = 1.03
$ws.Cells.Item(37, 11).Value = 9.25
$ws.Cells.Item(37, 13).Value = 4.55
$ws.Cells.Item(37, 14).Value = 1.5
$ws.Cells.Item(37, 15).Value = 2.4
$ws.Cells.Item(37, 16).Value = 1.28
$ws.Cells.Item(37, 17).Value = 3.35
$ws.Cells.Item(37, 18).Value = 1.44
$ws.Cells.Item(37, 19).Value = 2.6
$ws.Cells.Item(37, 20).Value = 14
$ws.Cells.Item(37, 21).Value = 18.5
$ws.Cells.Item(37, 22).Value = 10.25
$ws.Cells.Item(37, 23).Value = 35
$ws.Cells.Item(37, 24).Value = 19
$ws.Cells.Item(37, 26).Value = 9.25
$ws.Cells.Item(37, 27).Value = 7.5
$ws.Cells.Item(37, 28).Value = 10.75
$ws.Cells.Item(37, 29).Value = 32
$ws.Cells.Item(37, 30).Value = 175
$ws.Cells.Item(37, 31).Value = 11.75
$ws.Cells.Item(37, 32).Value = 14.5
$ws.Cells.Item(37, 33).Value = 9.25
$ws.Cells.Item(37, 34).Value = 25
$ws.Cells.Item(37, 35).Value = 16.5
$ws.Cells.Item(37, 36).Value = 19.5

# Row 48
$ws.Cells.Item(48, 7).Value = 8.5
$ws.Cells.Item(48, 14).Value = 1.6
$ws.Cells.Item(48, 15).Value = 2.05
$ws.Cells.Item(48, 18).Value = 2.13
$ws.Cells.Item(48, 19).Value = 1.64
$ws.Cells.Item(48, 20).Value = 17.5
$ws.Cells.Item(48, 23).Value = 175
$ws.Cells.Item(48, 24).Value = 80
$ws.Cells.Item(48, 25).Value = 70
$ws.Cells.Item(48, 26).Value = 12.5
$ws.Cells.Item(48, 27).Value = 8.5
$ws.Cells.Item(48, 28).Value = 19.5
$ws.Cells.Item(48, 29).Value = 90
$ws.Cells.Item(48, 31).Value = 6
$ws.Cells.Item(48, 32).Value = 5.2
$ws.Cells.Item(48, 34).Value = 6.4
$ws.Cells.Item(48, 36).Value = 24

# Row 49
$ws.Cells.Item(49, 7).Value = 4.9
$ws.Cells.Item(49, 8).Value = 3.6
$ws.Cells.Item(49, 9).Value = 1.6
$ws.Cells.Item(49, 12).Value = 1.31
$ws.Cells.Item(49, 13).Value = 3.15
$ws.Cells.Item(49, 14).Value = 1.8
$ws.Cells.Item(49, 15).Value = 1.8
$ws.Cells.Item(49, 16).Value = 1.37
$ws.Cells.Item(49, 17).Value = 2.5
$ws.Cells.Item(49, 18).Value = 2
$ws.Cells.Item(49, 19).Value = 1.73
$ws.Cells.Item(49, 20).Value = 10.75
$ws.Cells.Item(49, 21).Value = 23
$ws.Cells.Item(49, 22).Value = 13.5
$ws.Cells.Item(49, 23).Value = 65
$ws.Cells.Item(49, 24).Value = 40
$ws.Cells.Item(49, 25).Value = 40
$ws.Cells.Item(49, 27).Value = 6.2
$ws.Cells.Item(49, 28).Value = 13.5
$ws.Cells.Item(49, 29).Value = 60
$ws.Cells.Item(49, 30).Value = 400
$ws.Cells.Item(49, 31).Value = 5.8
$ws.Cells.Item(49, 32).Value = 6.4
$ws.Cells.Item(49, 33).Value = 6.8
$ws.Cells.Item(49, 34).Value = 9.75
$ws.Cells.Item(49, 35).Value = 10.5
$ws.Cells.Item(49, 36).Value = 21

# Row 51
$ws.Cells.Item(51, 7).Value = 2.13
$ws.Cells.Item(51, 8).Value = 3.05
$ws.Cells.Item(51, 9).Value = 3.3
$ws.Cells.Item(51, 10).Value = 1.11
$ws.Cells.Item(51, 11).Value = 5.6
$ws.Cells.Item(51, 12).Value = 1.58
$ws.Cells.Item(51, 13).Value = 2.24
$ws.Cells.Item(51, 14).Value = 2.7
$ws.Cells.Item(51, 15).Value = 1.4
$ws.Cells.Item(51, 16).Value = 1.64
$ws.Cells.Item(51, 17).Value = 2.15
$ws.Cells.Item(51, 18).Value = 2.33
$ws.Cells.Item(51, 19).Value = 1.54
$ws.Cells.Item(51, 20).Value = 4.1
$ws.Cells.Item(51, 21).Value = 6.8
$ws.Cells.Item(51, 22).Value = 8
$ws.Cells.Item(51, 23).Value = 17
$ws.Cells.Item(51, 24).Value = 20
$ws.Cells.Item(51, 25).Value = 45
$ws.Cells.Item(51, 26).Value = 4.8
$ws.Cells.Item(51, 27).Value = 5
$ws.Cells.Item(51, 28).Value = 19
$ws.Cells.Item(51, 29).Value = 101
$ws.Cells.Item(51, 30).Value = 101
$ws.Cells.Item(51, 31).Value = 5.4
$ws.Cells.Item(51, 32).Value = 12
$ws.Cells.Item(51, 33).Value = 11
$ws.Cells.Item(51, 34).Value = 40
$ws.Cells.Item(51, 35).Value = 35
$ws.Cells.Item(51, 36).Value = 60

# Row 54
$ws.Cells.Item(54, 10).Value = 1.06
$ws.Cells.Item(54, 11).Value = 10
$ws.Cells.Item(54, 18).Value = 1.8
$ws.Cells.Item(54, 19).Value = 1.91
$ws.Cells.Item(54, 20).Value = 8.5
$ws.Cells.Item(54, 30).Value = 251

# Row 62
$ws.Cells.Item(62, 7).Value = 3.3
$ws.Cells.Item(62, 8).Value = 3.4
$ws.Cells.Item(62, 9).Value = 2.02
$ws.Cells.Item(62, 12).Value = 1.23
$ws.Cells.Item(62, 13).Value = 3.35
$ws.Cells.Item(62, 14).Value = 1.7
$ws.Cells.Item(62, 15).Value = 1.91
$ws.Cells.Item(62, 20).Value = 11.75
$ws.Cells.Item(62, 21).Value = 19.5
$ws.Cells.Item(62, 22).Value = 11.25
$ws.Cells.Item(62, 26).Value = 11.5
$ws.Cells.Item(62, 27).Value = 6.7
$ws.Cells.Item(62, 31).Value = 8.25
$ws.Cells.Item(62, 32).Value = 10.5
$ws.Cells.Item(62, 34).Value = 19

# Row 63
$ws.Cells.Item(63, 8).Value = 4.35
$ws.Cells.Item(63, 13).Value = 3.6
$ws.Cells.Item(63, 14).Value = 1.62
$ws.Cells.Item(63, 15).Value = 2.02
$ws.Cells.Item(63, 18).Value = 1.8
$ws.Cells.Item(63, 20).Value = 7.5
$ws.Cells.Item(63, 21).Value = 7.1
$ws.Cells.Item(63, 24).Value = 11.5
$ws.Cells.Item(63, 26).Value = 13
$ws.Cells.Item(63, 29).Value = 80
$ws.Cells.Item(63, 31).Value = 17

# Row 64
$ws.Cells.Item(64, 8).Value = 3.25
$ws.Cells.Item(64, 9).Value = 2.6
$ws.Cells.Item(64, 11).Value = 8.5
$ws.Cells.Item(64, 28).Value = 17

# Row 65
$ws.Cells.Item(65, 7).Value = 1.7
$ws.Cells.Item(65, 8).Value = 3.8
$ws.Cells.Item(65, 14).Value = 1.95
$ws.Cells.Item(65, 15).Value = 1.9
$ws.Cells.Item(65, 16).Value = 1.36
$ws.Cells.Item(65, 17).Value = 3
$ws.Cells.Item(65, 26).Value = 11
$ws.Cells.Item(65, 31).Value = 13

# Row 71
$ws.Cells.Item(71, 11).Value = 15

# Row 72
$ws.Cells.Item(72, 7).Value = 1.91
$ws.Cells.Item(72, 8).Value = 3.6
$ws.Cells.Item(72, 9).Value = 3.55
$ws.Cells.Item(72, 10).Value = 1.05
$ws.Cells.Item(72, 11).Value = 8
$ws.Cells.Item(72, 12).Value = 1.24
$ws.Cells.Item(72, 13).Value = 3.65
$ws.Cells.Item(72, 14).Value = 1.72
$ws.Cells.Item(72, 15).Value = 2
$ws.Cells.Item(72, 16).Value = 1.35
$ws.Cells.Item(72, 17).Value = 2.92
$ws.Cells.Item(72, 18).Value = 1.65
$ws.Cells.Item(72, 19).Value = 2.1
$ws.Cells.Item(72, 20).Value = 8.5
$ws.Cells.Item(72, 22).Value = 8.25
$ws.Cells.Item(72, 23).Value = 17
$ws.Cells.Item(72, 24).Value = 14
$ws.Cells.Item(72, 25).Value = 23
$ws.Cells.Item(72, 26).Value = 8
$ws.Cells.Item(72, 27).Value = 7
$ws.Cells.Item(72, 28).Value = 13.5
$ws.Cells.Item(72, 29).Value = 55
$ws.Cells.Item(72, 30).Value = 350
$ws.Cells.Item(72, 31).Value = 11.75
$ws.Cells.Item(72, 32).Value = 20
$ws.Cells.Item(72, 33).Value = 12
$ws.Cells.Item(72, 34).Value = 50
$ws.Cells.Item(72, 35).Value = 30
$ws.Cells.Item(72, 36).Value = 35

# Row 73
$ws.Cells.Item(73, 7).Value = 1.7
$ws.Cells.Item(73, 8).Value = 3.9
$ws.Cells.Item(73, 9).Value = 4.3
$ws.Cells.Item(73, 10).Value = 1.04
$ws.Cells.Item(73, 11).Value = 8.5
$ws.Cells.Item(73, 12).Value = 1.21
$ws.Cells.Item(73, 13).Value = 3.9
$ws.Cells.Item(73, 14).Value = 1.65
$ws.Cells.Item(73, 15).Value = 2.12
$ws.Cells.Item(73, 16).Value = 1.32
$ws.Cells.Item(73, 17).Value = 3.1
$ws.Cells.Item(73, 18).Value = 1.65
$ws.Cells.Item(73, 19).Value = 2.1
$ws.Cells.Item(73, 20).Value = 8.5
$ws.Cells.Item(73, 21).Value = 9
$ws.Cells.Item(73, 22).Value = 8.25
$ws.Cells.Item(73, 23).Value = 13.5
$ws.Cells.Item(73, 24).Value = 12
$ws.Cells.Item(73, 25).Value = 22
$ws.Cells.Item(73, 26).Value = 8.5
$ws.Cells.Item(73, 27).Value = 7.7
$ws.Cells.Item(73, 28).Value = 14
$ws.Cells.Item(73, 29).Value = 55
$ws.Cells.Item(73, 30).Value = 350
$ws.Cells.Item(73, 31).Value = 14.5
$ws.Cells.Item(73, 32).Value = 26
$ws.Cells.Item(73, 33).Value = 14
$ws.Cells.Item(73, 34).Value = 70
$ws.Cells.Item(73, 35).Value = 37
$ws.Cells.Item(73, 36).Value = 37
